$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add hours logged for the two new dates (rows 15 and 16, column D)
$ws.Range("D15").Value = 6
$ws.Range("D16").Value = 7

# Update the active cell selection to match the saved view state
$ws.Range("I14").Select()
